# Refresh the cryptocurrency price/volume table with newly scraped values.
# A few coins also changed rank (rows re-sorted), which shows up as the
# Coin/Link/Price/Volume cells for those rows changing together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cell, new literal text, whether the text looks like a plain number).
# The sheet stores prices as text (e.g. "1.000", "30.360.44") so any value
# that Excel would otherwise auto-convert to a number needs the cell forced
# to text first; we restore the original cell style right after so no
# formatting changes are introduced beyond the requested value edits.
$updates = @(
    @{ Cell = 'D2'; Value = '30.360.44'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.863.70'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -0.48%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  +0.02%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '234.68'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -1.45%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '1.000'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +0.04%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.4776'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -0.16%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.2775'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -1.65%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.06506'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -0.02%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '1.863.32'; ForceText = $false },
    @{ Cell = 'E10'; Value = '  -0.55%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.07436'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -0.25%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  -2.09%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '5.032'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -1.43%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '86.61'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -1.74%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '0.6368'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -2.64%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '30.339.59'; ForceText = $false },
    @{ Cell = 'E16'; Value = '  +0.01%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '1.000'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'B18'; Value = 'Avalanche'; ForceText = $false },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; ForceText = $false },
    @{ Cell = 'D18'; Value = '12.93'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -2.86%  '; ForceText = $false },
    @{ Cell = 'B19'; Value = 'BitcoinCash'; ForceText = $false },
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; ForceText = $false },
    @{ Cell = 'D19'; Value = '233.36'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  +6.48%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '0.000007413'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -2.26%  '; ForceText = $false },
    @{ Cell = 'B21'; Value = 'BinanceUSD'; ForceText = $false },
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText = $false },
    @{ Cell = 'D21'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  +0.15%  '; ForceText = $false },
    @{ Cell = 'B22'; Value = 'Uniswap'; ForceText = $false },
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; ForceText = $false },
    @{ Cell = 'D22'; Value = '5.110'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -3.53%  '; ForceText = $false },
    @{ Cell = 'B23'; Value = 'Chainlink'; ForceText = $false },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; ForceText = $false },
    @{ Cell = 'D23'; Value = '6.119'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -1.11%  '; ForceText = $false },
    @{ Cell = 'B24'; Value = 'Monero'; ForceText = $false },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false },
    @{ Cell = 'D24'; Value = '168.23'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  +0.30%  '; ForceText = $false },
    @{ Cell = 'B25'; Value = 'Cosmos'; ForceText = $false },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; ForceText = $false },
    @{ Cell = 'D25'; Value = '9.282'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -0.51%  '; ForceText = $false },
    @{ Cell = 'B26'; Value = 'EthereumClassic'; ForceText = $false },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText = $false },
    @{ Cell = 'D26'; Value = '18.09'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -1.82%  '; ForceText = $false },
    @{ Cell = 'B27'; Value = 'LidoDAOToken'; ForceText = $false },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; ForceText = $false },
    @{ Cell = 'D27'; Value = '1.888'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -4.35%  '; ForceText = $false },
    @{ Cell = 'B28'; Value = 'Stellar'; ForceText = $false },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false },
    @{ Cell = 'D28'; Value = '0.1046'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +11.76%  '; ForceText = $false },
    @{ Cell = 'B29'; Value = 'Toncoin'; ForceText = $false },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false },
    @{ Cell = 'D29'; Value = '1.382'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -4.41%  '; ForceText = $false },
    @{ Cell = 'B30'; Value = 'InternetComputer(DFINITY)'; ForceText = $false },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false },
    @{ Cell = 'D30'; Value = '4.245'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -1.62%  '; ForceText = $false },
    @{ Cell = 'B31'; Value = 'Filecoin'; ForceText = $false },
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false },
    @{ Cell = 'D31'; Value = '3.946'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -2.08%  '; ForceText = $false },
    @{ Cell = 'B32'; Value = 'Hedera'; ForceText = $false },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false },
    @{ Cell = 'D32'; Value = '0.04961'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -2.25%  '; ForceText = $false },
    @{ Cell = 'B33'; Value = 'ARBITRUM'; ForceText = $false },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false },
    @{ Cell = 'D33'; Value = '1.161'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -3.71%  '; ForceText = $false },
    @{ Cell = 'B34'; Value = 'ImmutableX'; ForceText = $false },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false },
    @{ Cell = 'D34'; Value = '0.7339'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -2.50%  '; ForceText = $false },
    @{ Cell = 'B35'; Value = 'Frax'; ForceText = $false },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; ForceText = $false },
    @{ Cell = 'D35'; Value = '0.9996'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  +0.15%  '; ForceText = $false },
    @{ Cell = 'B36'; Value = 'HuobiToken'; ForceText = $false },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; ForceText = $false },
    @{ Cell = 'D36'; Value = '2.712'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'B37'; Value = 'VeChain'; ForceText = $false },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false },
    @{ Cell = 'D37'; Value = '0.01925'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +5.47%  '; ForceText = $false },
    @{ Cell = 'B38'; Value = 'MXToken'; ForceText = $false },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; ForceText = $false },
    @{ Cell = 'D38'; Value = '2.635'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +0.88%  '; ForceText = $false },
    @{ Cell = 'B39'; Value = 'TrustWalletToken'; ForceText = $false },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.9120'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +0.60%  '; ForceText = $false },
    @{ Cell = 'B40'; Value = 'RenderToken'; ForceText = $false },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false },
    @{ Cell = 'D40'; Value = '2.016'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -2.53%  '; ForceText = $false },
    @{ Cell = 'B41'; Value = 'Quant'; ForceText = $false },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false },
    @{ Cell = 'D41'; Value = '105.61'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -1.14%  '; ForceText = $false },
    @{ Cell = 'B42'; Value = 'PaxDollar'; ForceText = $false },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; ForceText = $false },
    @{ Cell = 'D42'; Value = '0.9958'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -0.78%  '; ForceText = $false },
    @{ Cell = 'B43'; Value = 'TheSandbox'; ForceText = $false },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.4164'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -2.51%  '; ForceText = $false },
    @{ Cell = 'B44'; Value = 'FraxShare'; ForceText = $false },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; ForceText = $false },
    @{ Cell = 'D44'; Value = '5.560'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -5.74%  '; ForceText = $false },
    @{ Cell = 'B45'; Value = 'Aptos'; ForceText = $false },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false },
    @{ Cell = 'D45'; Value = '7.102'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -3.71%  '; ForceText = $false },
    @{ Cell = 'B46'; Value = 'Aave'; ForceText = $false },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText = $false },
    @{ Cell = 'D46'; Value = '61.35'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -4.33%  '; ForceText = $false },
    @{ Cell = 'B47'; Value = 'EnergySwap'; ForceText = $false },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false },
    @{ Cell = 'D47'; Value = '8.866'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -0.76%  '; ForceText = $false },
    @{ Cell = 'B48'; Value = 'Algorand'; ForceText = $false },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false },
    @{ Cell = 'D48'; Value = '0.1214'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -5.45%  '; ForceText = $false },
    @{ Cell = 'B49'; Value = 'NEARProtocol'; ForceText = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText = $false },
    @{ Cell = 'D49'; Value = '1.410'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -3.95%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'Elrond'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'; ForceText = $false },
    @{ Cell = 'D50'; Value = '33.36'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -0.67%  '; ForceText = $false },
    @{ Cell = 'B51'; Value = 'Cronos'; ForceText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.05626'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -0.01%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $u.Value
    }
}
